$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update price column (D) values for rows 27-34 (correct an over-charge bug)
$ws.Range("D27").Value = 5215.875
$ws.Range("D28").Value = 6520.5
$ws.Range("D29").Value = 5588.625
$ws.Range("D30").Value = 7079.625
$ws.Range("D31").Value = 5961.375
$ws.Range("D32").Value = 7637.438
$ws.Range("D33").Value = 6520.5
$ws.Range("D34").Value = 8759.625
